$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FilesTab (row 4) query text was reordered: Format / File Type / Size columns
# moved earlier (right after File Name), matching the "fix case files
# comparison icdc" update referenced in the commit message.
$newFilesTabQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Poodle'] 
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '') AS `File Name`,
       coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
'@

$ws.Range("B4").Value2 = $newFilesTabQuery

# Keep row 4 at Excel's capped max row height (matches the other
# long-query rows, which are already pinned at 409.5pt) instead of
# whatever uncapped auto-fit height the new text would otherwise produce.
$ws.Rows.Item(4).RowHeight = 409.5

# Reflect the author's final cursor position (scrolled down, landing on the
# FilesTab query cell that was just edited).
$ws.Range("B4").Select()
